$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$tb = $s.Shapes.AddTextbox(1, 59.338582677165356, 238.65551181102362, 384.4251968503937, 42.425196850393704)
$tf = $tb.TextFrame
$tf.TextRange.Text = "Doceree Karasuno Zip file link"
$tf.WordWrap = $true
$tf.AutoSize = 1
$tf.MarginLeft = 7.2
$tf.MarginRight = 7.2
$tf.MarginTop = 7.2
$tf.MarginBottom = 7.2
